$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = "2025-04-28 23:37:52"
$ws.Range("B33").Value = 178
